## BAU Fraction of Heat from CHP.xlsx - apply commit "Dropped in WRI's input data
## for first draft of 2.1. made several fixes to model run errors".
##
## Summary of edits:
##  1. "About" sheet: rewrite the Notes paragraph (split "We do not have data..."
##     into two strings, capitalize "So we assume..."), re-insert the unchanged
##     trailing notes, and add a new closing note row.
##  2. "BFoHfC" sheet: label the (previously blank) header cell A1, make it bold
##     + wrap, widen column A and heighten row 1 to fit; then add four new fuel
##     rows (crude oil, heavy or residual fuel oil, LPG propane or butane,
##     hydrogen) that mirror the existing rows' "=Data!$A$6" formula across all
##     year columns (B:AK).
##
## NOTE on write order: the new shared strings must be *created* (i.e. the cell
## that first uses each string must be written) in the same order they end up
## in the final workbook, since the string table is rebuilt in first-use order
## on save. The target order is: crude oil, heavy or residual fuel oil, LPG
## propane or butane, hydrogen, then the two split "We do not have..." notes
## sentences, then "So we assume...", then the new BFoHfC header text -- so the
## BFoHfC fuel rows are written first, the About notes second, and the BFoHfC
## header text last.

$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsBFoHfC = $wb.Worksheets.Item("BFoHfC")

## --- list of year columns used across BFoHfC data rows (B..AK) ---
$yearCols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z", `
              "AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK")

## ------------------------------------------------------------------
## 1. BFoHfC: four new fuel-type rows (8-11), each = Data!$A$6 across B:AK
## ------------------------------------------------------------------
$newFuels = @("crude oil", "heavy or residual fuel oil", "LPG propane or butane", "hydrogen")

$row = 8
foreach ($fuel in $newFuels) {
    $wsBFoHfC.Range("A$row").Value = $fuel

    foreach ($col in $yearCols) {
        $wsBFoHfC.Range("$col$row").Formula = "=Data!`$A`$6"
    }
    # The formulas above point at a single-cell reference whose own number
    # format is "0.000" -- Excel inherits that format onto a blank/General
    # destination cell the first time such a formula is entered. The source
    # workbook keeps these new cells on the default General format, so clear
    # the inherited number formatting back off again.
    $wsBFoHfC.Range("B$row`:AK$row").ClearFormats()

    $row = $row + 1
}

## ------------------------------------------------------------------
## 2. About sheet: rewrite the Notes paragraph text (rows 10-12), shift the
##    unchanged trailing sentences down (rows 13-14), and append the final
##    sentence on a new row 15.
## ------------------------------------------------------------------
$wsAbout.Range("A10").Value = "We do not have data on how Combined Heat and Power (CHP)"
$wsAbout.Range("A11").Value = " fraction varies by fuel type."
$wsAbout.Range("A12").Value = "So we assume the overall CHP fraction applies to all fuel types."
$wsAbout.Range("A13").Value = "It does not matter if some of these fuel types are not used"
$wsAbout.Range("A14").Value = "for district heat; the fraction of each fuel that is used is"
$wsAbout.Range("A15").Value = "governed by another input variable."

## Leave the cursor parked on the last edited note row, matching the saved
## selection state in the target workbook.
$wsAbout.Activate()
$wsAbout.Range("A13").Select()

## ------------------------------------------------------------------
## 3. BFoHfC: label + style the header cell A1, widen column A, heighten row 1
## ------------------------------------------------------------------
$wsBFoHfC.Range("A1").Value = "Fraction of Heat from CHP (dimensionless)"
$wsBFoHfC.Range("A1").Font.Bold = $true
$wsBFoHfC.Range("A1").WrapText = $true

$wsBFoHfC.Columns.Item(1).ColumnWidth = 24.5
$wsBFoHfC.Rows.Item(1).RowHeight = 30

Write-Host "Applied BAU Fraction of Heat from CHP edits"
